$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 16 (Case 5471 / QUESADA 3212), shifting all subsequent rows up
$ws.Rows.Item(16).Delete()
